$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.922.36"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.419.33"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "551.86"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "137.56"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +3.12%  "
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").Value = "5.68"
$ws.Range("E10").Value = "  -2.96%  "
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "25.52"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").Value = "2.847.08"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "59.838.23"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "2.422.10"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "4.41"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").Value = "  -4.97%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "66.45"
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("D24").Value = "0.173"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").Value = "8.67"
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "1.38"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").Value = "0.0₃0777"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").Value = "169.22"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "18.70"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "4.23"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "314.44"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "0.409"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "19.60"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "0.387"
$ws.Range("E48").Value = "  -5.80%  "
$ws.Range("D49").Value = "17.70"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "11.05"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -0.21%  "
